$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "MEC-2B-Maq. Term. FLuxo" class block from column D (rows 3-4)
# up/over to column E (rows 2-3), clearing the vacated cells with "-".
$ws.Range("E2").Value = "MEC-2B-Maq. Term. FLuxo"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "MEC-2B-Maq. Term. FLuxo"
$ws.Range("D4").Value = "-"
